$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, border, centered) from AC1
# onto the new header cells so they reuse the same style (no new style
# entries created).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the win/loss/tie record for every data row (2-66)
$ws.Range("AD2:AD66").Value = 52
$ws.Range("AE2:AE66").Value = 110
$ws.Range("AF2:AF66").Value = 0
